$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the "DATE" label that used to sit in A1 (header row).
#    The header row now reads: (blank) | DAILY ROUTINE | DOCUMENT | TECHNICAL
# ------------------------------------------------------------------
$ws.Range("A1").ClearContents()

# ------------------------------------------------------------------
# 2. Append the new "daily routine" rows (38-50) below the existing data.
#    Column A holds dates formatted the same way as the rest of the sheet
#    (copy the date style from A2 so the same cellXf / numFmt is reused).
#    Column B holds the corresponding activity text.
# ------------------------------------------------------------------

# Helper: copy the date-formatted style used throughout column A
$ws.Range("A2").Copy() | Out-Null

$dateRows = 38..50
foreach ($r in $dateRows) {
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# -- Column A first (matches the order the original author typed the data
#    in: all the dates, then the description text in column B) --

# Rows 38-41: plain dates (Excel recognises them as valid dates)
$ws.Range("A38").Value = 44014
$ws.Range("A39").Value = 44106
$ws.Range("A40").Value = 44137
$ws.Range("A41").Value = 44167

# Rows 42-49: dates typed day/month/year with a day > 12, so Excel cannot
# parse them as dates and stores them as plain text (keeping the date
# number format already applied to the cell).
$ws.Range("A42").Value = "13/2/2020"
$ws.Range("A43").Value = "14/2/2020"
$ws.Range("A44").Value = "17/2/2020"
$ws.Range("A45").Value = "18/2/2020"
$ws.Range("A46").Value = "19/2/2020"
$ws.Range("A47").Value = "20/2/2020"
$ws.Range("A48").Value = "21/2/2020"
$ws.Range("A49").Value = "24/2/2020"

# Row 50: left blank, only keeps the date style that was pasted in above.

# -- Column B second --
$ws.Range("B39").Value = "training"
$ws.Range("B40").Value = "training"
$ws.Range("B41").Value = "training"
$ws.Range("B42").Value = "training"
$ws.Range("B43").Value = "training"
$ws.Range("B44").Value = "training"
$ws.Range("B45").Value = "training"
$ws.Range("B46").Value = "training"
$ws.Range("B47").Value = "training"
$ws.Range("B48").Value = "project connectivity"
$ws.Range("B49").Value = "java script"

# ------------------------------------------------------------------
# 3. Update the active selection to match the end state (B49).
# ------------------------------------------------------------------
$ws.Range("B49").Select()
